$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 177, shifting rows 177:264 down to 178:265
$ws.Rows.Item(177).Insert()

# Fill in the new row 177 with the new data
$ws.Range("A177").Value = 7
$ws.Range("B177").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C177").Value = "Ñuble"
$ws.Range("D177").Value = 44572
$ws.Range("E177").Value = 16
$ws.Range("F177").Value = 100114001
$ws.Range("G177").Value = "Papa"
$ws.Range("H177").Value = "Asterix"
$ws.Range("I177").Value = "1a nueva(o)"
$ws.Range("J177").Value = 240
$ws.Range("K177").Value = 7500
$ws.Range("L177").Value = 8000
$ws.Range("M177").Value = 7750
$ws.Range("N177").Value = "`$/saco 25 kilos"
$ws.Range("O177").Value = "Región del Maule"
$ws.Range("P177").Value = 310
$ws.Range("Q177").Value = 25
$ws.Range("R177").Value = "Hortaliza"

# Ensure the date style (numFmtId 165) for D177 matches other D column cells
$ws.Range("D177").NumberFormat = $ws.Range("D178").NumberFormat
